$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G (date serial) for every data row 2..274 from 45554 to 45555
for ($i = 2; $i -le 274; $i++) {
    $ws.Cells.Item($i, 7).Value = 45555
}

# Update column E and H pairs that changed value (balance updates)
$changes = @{
    5 = 1346.85
    8 = 35172.26
    15 = 1049.71
    17 = 1037.45
    18 = 7.54
    19 = 5.78
    24 = 175.86
    43 = 458.95
    57 = 2908.32
    58 = 18.88
    59 = 2.8
    60 = 1080
    99 = 1556.8
    103 = 1070.1199999999999
    104 = 63.67
    108 = 2971.5
    132 = 985.85
    143 = 1209.06
    158 = 84.06
    165 = 5227.63
    173 = 2095.61
    198 = 19.07
    220 = 4.38
    226 = 202.69
    231 = 753.02
    235 = 536.77
    255 = 27060.240000000002
    264 = 2804.08
    265 = 1885.57
    270 = 22.84
    271 = 1244.58
    273 = 1576.08
    274 = 11.56
}
foreach ($row in $changes.Keys) {
    $val = $changes[$row]
    $ws.Cells.Item($row, 5).Value = $val
    $ws.Cells.Item($row, 8).Value = $val
}

# Rename sheet to reflect new export timestamp
$ws.Name = "IClientBalance-20240920-092250-"
